$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("MAIN_CONTROLLER")
$ws2 = $wb.Worksheets.Item("DATASHEET")

# --- DATASHEET: update FOS5.xlsx -> FOS3.xlsx (D3) ---
$ws2.Range("D3").Value = "FOS3.xlsx"

# --- MAIN_CONTROLLER: add new "ApplicationType" column (F) ---
$ws1.Range("F1").Value = "ApplicationType"

# New row 4 (copy of row 3 data) with ApplicationType = WEB.
# Filled before the MOBILE cells below so the shared-string table
# ends up with WEB before MOBILE, matching the source order.
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "N"
$ws1.Range("C4").Value = "local"
$ws1.Range("D4").Value = "Calculator"
$ws1.Range("E4").Value = "Calculator"
$ws1.Range("F4").Value = "WEB"

# Existing rows 2 & 3 get ApplicationType = MOBILE
$ws1.Range("F2").Value = "MOBILE"
$ws1.Range("F3").Value = "MOBILE"

# Columns E and F auto-fit to their (now wider) content, same as Excel
# does automatically when new data is entered.
$ws1.Columns(5).ColumnWidth = 8.42
$ws1.Columns(6).ColumnWidth = 13.25

# --- Selections / active sheet ---
# DATASHEET is no longer the active tab; its selection moves to E10.
$ws2.Range("E10").Select()

# MAIN_CONTROLLER becomes the active tab with F2 selected.
$ws1.Activate()
$ws1.Range("F2").Select()
